$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Customer Class" section (rows 3-6): grading completed, award full points
# in column E (Total Points) to match column D (Points for grading).
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# "Product Class" section (rows 10-14): grading completed, award full
# points in column E to match column D.
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the final selection on E15, matching the saved view state.
$ws.Range("E15").Select()
